# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the
# ec2463c9-00b0-4470-8a79-73491930967a row (row 4) on both
# the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-20 18:33:32"
$wsZhCn.Range("H4").Value = "2016-03-20 18:33:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-20 18:33:35"
$wsDeDe.Range("H4").Value = "2016-03-20 18:34:00"
